$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Original structure (5 paragraphs):
#   P1: "Uma Rádio necessita de um software..."
#   P2: (empty)
#   P3: "O Sistema deve reproduzir músicas..."
#   P4: "Há necessidade de Módulos bem definidos..."
#   P5: "O sistema deve ser capaz de funcionar..." + bookmark (_GoBack)
#
# Target structure (14 paragraphs) is built below, reusing P1/P2/P3/P4/P5
# where possible (to preserve paragraph identity/bookmark) and inserting
# new paragraphs for everything else.
# ---------------------------------------------------------------------

# --- T1: replace text of paragraph 1 ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "Este projeto visa o desenvolvimento de um software totalmente nacional com a finalidade de atender os requisitos básicos de uma rádio comunitária."

# --- T2: paragraph 2 stays empty, untouched ---
$p2 = $d.Paragraphs.Item(2)

# --- T3: insert new paragraph (text) after paragraph 2 ---
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "Há mais de 4 mil rádios comunitárias no Brasil, muitas das quais utilizam software para automatizar tarefas rotineiras, dentre elas reproduzir arquivos de áudio, inserção de vinhetas e spots, agendamento de eventos automatizados para os momentos em que não necessitar de intervenção humana."

# --- T4: insert new empty paragraph after T3 (no run) ---
$r3 = $p3.Range
$r3.Collapse(0)
$r3.Text = "`r"

# --- T5: the old "O Sistema deve reproduzir..." paragraph is now index 5. Replace its text. ---
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "Uma Rádio necessita de um sistema para contemplar operações rotineiras, de fácil manuseio, intuitivo a ponto de qualquer pessoa ser capaz de realizar as operações em segundos sem necessidade de prévio treinamento. O operador de áudio deve encontrar facilmente as diversas"

# --- T6: insert new paragraph (continuation text) after T5 ---
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "funcionalidades presentes no software, sendo organizado em módulos que podem interagir entre si mediante permissão previamente configurada."

# --- T7: insert new empty paragraph after T6 (no run) ---
$r6 = $p6.Range
$r6.Collapse(0)
$r6.Text = "`r"

# --- T8: insert new paragraph (text) after T7 ---
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "O sistema deve possuir um módulo com player de áudio padrão onde reproduzirá arquivos de áudio que estarão listadas em uma playlist, este será o mais utilizado e merece destaque."

# --- T9: insert new empty paragraph after T8 (no run) ---
$r8 = $p8.Range
$r8.Collapse(0)
$r8.Text = "`r"

# --- T10: the old "Há necessidade de Módulos..." paragraph is now index 10. Replace its text. ---
$p10 = $d.Paragraphs.Item(10)
$p10.Range.Text = "Há necessidade de Módulos bem definidos, para a inserção das Vinhetas, Spots e Programações gravadas, mais ainda gerenciamento dos eventos automatizados, gerenciamento das vinhetas, spots e programações gravadas, possibilitando reproduzi-los de maneira eficiente, redirecionando para a Playlist ao término. Todos os Módulos precisam ter a possibilidade de serem agendados e executados pelo evento automatizado."

# --- T11: insert new empty paragraph after T10 (no run) ---
$r10 = $p10.Range
$r10.Collapse(0)
$r10.Text = "`r"

# --- T12: insert new paragraph (text, vertAlign baseline) after T11 ---
$p11 = $d.Paragraphs.Item(11)
$p11.Range.InsertParagraphAfter()
$p12 = $d.Paragraphs.Item(12)
$p12.Range.Text = "O sistema deve ser capaz de agendar eventos complexos, execuções e finalizações de arquivos e sistemas de forma automatizada, mais ainda agendar o desligamento automático do computador."
$p12.Range.Font.Subscript = $false

# --- T13: the old final paragraph ("O sistema deve ser capaz de funcionar..." + bookmark) is now
#     index 13. Clear its run/text but keep the paragraph (and bookmark) intact. ---
$p13 = $d.Paragraphs.Item(13)
$s = $p13.Range.Start
$e = $p13.Range.End
$clearRange = $d.Range($s, $e - 1)
$clearRange.Delete()

# --- T14: insert new paragraph (final text) after T13, i.e. after the bookmark ---
$p13b = $d.Paragraphs.Item(13)
$p13b.Range.InsertParagraphAfter()
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Text = "O sistema deve ser capaz de funcionar em outras plataformas, como Windows e Linux, dando maior liberdade para a escolha do Sistema Operacional, deve funcionar 24 horas por dia, podendo funcionar em segundo plano, apresentando notificação caso seja finalizado."
